# "Modificacion en reporte de stock para ver marca"
# - Row 43 (task "reporte movimiento cliente - los creditos van en haber")
#   gets a responsible person (Agustina) and status (en proceso), same as
#   row 44 already had.
# - Row 44 (task "en stock agregar codigo - descripcion") status switches
#   from the text "en proceso" to a completed (100%) percentage value,
#   reusing the same percentage number format as the other finished rows.
# - Column A is widened (no longer relying on best-fit) to make room to
#   show the "marca" (brand) text.
# - The active selection moves up one row to C44.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Widen column A (remove best-fit, make it a bit wider / custom width).
$ws.Columns.Item(1).ColumnWidth = 76.67

# Row 43: add Responsable + Estado, matching the existing "en proceso" entries.
$ws.Range("B43").Value = "Agustina"
$ws.Range("C43").Value = "en proceso"

# Row 44: mark as complete (100%) instead of "en proceso".
$ws.Range("C44").NumberFormat = "0%"
$ws.Range("C44").Value = 1

# Update the selected / active cell to C44.
$ws.Range("C44").Select()
